# Weekly update: insert a new record at row 13 (pushing all subsequent
# weekly records down by one row) and keep the previously-last record
# (old row 143) as the new last row (144).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13; everything from old row 13
# downward (through old row 143) shifts down to rows 14..144.
$ws.Rows("13:13").Insert()

# Populate the new row 13 with this week's record (same dimensions /
# metadata as the rest of the "Cebollín" table, new Fecha + Volumen).
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44552
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100112037
$ws.Range("G13").Value = "Cebollín"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 900
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 950
$ws.Range("N13").Value = "$/paquete 6 unidades"
$ws.Range("O13").Value = "Provincia del Elquí"
$ws.Range("P13").Value = 158
$ws.Range("Q13").Value = 6
$ws.Range("R13").Value = "Hortaliza"
